# fix: fix rubics function
#
# Row 13 used to hold the (buggy) "scrambled_rubics_image.jpg" metrics row and
# row 14 held "scrambled_image_prime.jpg". The rubics computation was fixed, so
# a proper "rubics.jpg" row now takes row 13 (with freshly computed metrics),
# and the old row-13 data (scrambled_image_prime.jpg) shifts down to row 14,
# dropping the old row 14 content entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: rubics.jpg with freshly computed metrics
$ws.Range("A13").Value = "rubics.jpg"
$ws.Range("B13").Value = 8.4218070881881353
$ws.Range("C13").Value = 7.4981904775947745
$ws.Range("D13").Value = 1.1231786006708158
$ws.Range("E13").Value = 0.1756744384765625

# Row 14 now carries what used to be row 13's data (scrambled_image_prime.jpg)
$ws.Range("A14").Value = "scrambled_image_prime.jpg"
$ws.Range("B14").Value = 566.76161579289055
$ws.Range("C14").Value = 7.3247438424006228
$ws.Range("D14").Value = 77.376305299864143
$ws.Range("E14").Value = 0.35858154296875

# Column A auto-shrinks slightly now that the longest label is
# "scrambled_image_prime.jpg" (25 chars) instead of
# "scrambled_rubics_image.jpg" (26 chars).
$ws.Columns.Item(1).ColumnWidth = 18.5
